# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" row at the top of the "总计" (summary) sheet's
#     data, shifting the existing quarters down by one.
#  2. Insert a new "2022-Q4" worksheet (with the Q4 fund-holdings detail)
#     right after "总计" / before "2022-Q3", shifting every later sheet
#     down by one tab position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert row 2 for "2022-Q4"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

# Column A carries the bold/bordered "index" style (s="2") - copy it from
# the row below (still formatted) onto the freshly inserted row.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value2 = 0
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 9
$summary.Cells.Item(2, 4).Value2 = 0.29

# Renumber the index column for the rows that got pushed down (used to be
# 0..5 starting at row 2, now 1..6 starting at row 3).
for ($r = 3; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value2 = $r - 2
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, positioned before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)

# Duplicate the "2022-Q3" sheet (keeps headers/column styles identical)
# and drop the copy in front of it, then rename.
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The Q3 sheet only has 6 data rows (rows 2-7); Q4 needs 9 (rows 2-10), so
# insert 3 more rows at the bottom, matching the existing row format.
$q4.Rows.Item(8).Insert()
$q4.Rows.Item(8).ClearFormats()
$q4.Rows.Item(9).Insert()
$q4.Rows.Item(9).ClearFormats()
$q4.Rows.Item(10).Insert()
$q4.Rows.Item(10).ClearFormats()

$q4.Range("A7").Copy()
$q4.Range("A8:A10").PasteSpecial(-4122)

# ---- header row (unchanged from Q3, rewritten for clarity/safety) -----
$q4.Cells.Item(1, 2).Value2 = "基金代码"
$q4.Cells.Item(1, 3).Value2 = "基金名称"
$q4.Cells.Item(1, 4).Value2 = "基金规模"
$q4.Cells.Item(1, 5).Value2 = "股票总仓位"
$q4.Cells.Item(1, 6).Value2 = "仓位占比"
$q4.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value2 = "仓位排名"

# ---- data rows ----------------------------------------------------
# Columns D, E, F, G are stored as *text* (not numbers) in the source
# workbook, so every value is assigned with a leading apostrophe to force
# text storage, then ClearFormats() strips the incidental "@" number
# format that the apostrophe trick leaves behind.
$q4Data = @(
    @(0, "005014", "泰康景泰回报混合A",             "8.86", "32.91", "1.34", "0.1187", 6),
    @(1, "519991", "长信双利优选混合A",             "1.11", "91.65", "4.19", "0.0465", 9),
    @(2, "006396", "长信双利优选混合E",             "1.11", "91.65", "4.19", "0.0465", 9),
    @(3, "512040", "富国中证价值ETF",               "3.39", "99.29", "1.08", "0.0366", 5),
    @(4, "009937", "东方欣益一年持有期偏债混合A",   "1.93", "22.19", "0.73", "0.0141", 7),
    @(5, "007142", "嘉合稳健增长灵活配置混合C",     "0.39", "88.99", "2.76", "0.0108", 6),
    @(6, "007141", "嘉合稳健增长灵活配置混合A",     "0.26", "88.99", "2.76", "0.0072", 6),
    @(7, "005015", "泰康景泰回报混合C",             "0.37", "32.91", "1.34", "0.0050", 6),
    @(8, "009938", "东方欣益一年持有期偏债混合C",   "0.31", "22.19", "0.73", "0.0023", 7)
)

$row = 2
foreach ($rec in $q4Data) {
    $q4.Cells.Item($row, 1).Value2 = $rec[0]
    $q4.Cells.Item($row, 2).Value2 = "'" + $rec[1]
    $q4.Cells.Item($row, 2).ClearFormats()
    $q4.Cells.Item($row, 3).Value2 = $rec[2]
    $q4.Cells.Item($row, 4).Value2 = "'" + $rec[3]
    $q4.Cells.Item($row, 4).ClearFormats()
    $q4.Cells.Item($row, 5).Value2 = "'" + $rec[4]
    $q4.Cells.Item($row, 5).ClearFormats()
    $q4.Cells.Item($row, 6).Value2 = "'" + $rec[5]
    $q4.Cells.Item($row, 6).ClearFormats()
    $q4.Cells.Item($row, 7).Value2 = "'" + $rec[6]
    $q4.Cells.Item($row, 7).ClearFormats()
    $q4.Cells.Item($row, 8).Value2 = $rec[7]
    $row++
}

# Restore column A's index style (ClearFormats on B/D/E/F/G calls above
# don't touch column A, but re-apply defensively in case row 2 lost it).
$q4.Range("A7").Copy()
$q4.Range("A2:A10").PasteSpecial(-4122)
